# Trade #14 (MarketMaking, 2026-02-17 13:10:15) closed out.
# Update the Summary, Strategy Status, All Trades, and MarketMaking sheets
# to reflect the closed trade's final P&L / capital figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.37   # Current Capital
$summary.Range("B4").Value = -0.63     # Total P&L $
$summary.Range("B5").Value = -0.9      # Total P&L %
$summary.Range("B6").Value = 14        # Total Trades
$summary.Range("B8").Value = 9         # Losing Trades
$summary.Range("B9").Value = 28.57     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.37      # Capital
$status.Range("D4").Value = 14         # Trades
$status.Range("E4").Value = -0.63      # P&L $
$status.Range("F4").Value = -0.63      # P&L %
$status.Range("G4").Value = 28.57      # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - Trade #14 row (row 15)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G15").Value = 0.3            # Exit Price
$allTrades.Range("H15").Value = "CLOSED"       # Status
$allTrades.Range("I15").Value = -43.3962       # P&L %
$allTrades.Range("J15").Value = -0.23          # P&L $
$allTrades.Range("K15").Value = 99.37          # Capital After
$allTrades.Range("P15").Value = "early_exit"   # Exit Reason
$allTrades.Range("Q15").Value = 7.16           # Duration (min)

# ---------------------------------------------------------------------
# MarketMaking sheet - Trade #14 row (row 15) - mirrors All Trades
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G15").Value = 0.3            # Exit Price
$mm.Range("H15").Value = "CLOSED"       # Status
$mm.Range("I15").Value = -43.3962       # P&L %
$mm.Range("J15").Value = -0.23          # P&L $
$mm.Range("K15").Value = 99.37          # Capital After
$mm.Range("P15").Value = "early_exit"   # Exit Reason
$mm.Range("Q15").Value = 7.16           # Duration (min)
